$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161 (shifts existing rows 161-286 down to 162-287)
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with its data
$ws.Range("A161").Value = 11
$ws.Range("B161").Value = "Vega Monumental Concepción"
$ws.Range("C161").Value = "Bíobío"
$ws.Range("D161").Value = 45072
$ws.Range("E161").Value = 8
$ws.Range("F161").Value = 100112003
$ws.Range("G161").Value = "Ajo"
$ws.Range("H161").Value = "Chino"
$ws.Range("I161").Value = "1a (cosecha)"
$ws.Range("J161").Value = 150
$ws.Range("K161").Value = 14000
$ws.Range("L161").Value = 15000
$ws.Range("M161").Value = 14467
$ws.Range("N161").Value = "`$/caja 10 kilos"
$ws.Range("O161").Value = "China"
$ws.Range("P161").Value = 1447
$ws.Range("Q161").Value = 10
$ws.Range("R161").Value = "Hortaliza"
